$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J7").Value = "sdgd fad vdfdf"
$ws.Range("K7").Value = "asdfas"
$ws.Range("L7").Value = "df"
$ws.Range("M7").Value = "fas"
$ws.Range("N7").Value = "dfasd"
$ws.Range("O7").Value = "fasd"

$ws.Range("J8").Value = "sdfas"
$ws.Range("O8").Value = "fasd"

$ws.Range("J9").Value = "df"
$ws.Range("O9").Value = "asdf"

$ws.Range("J10").Value = "asdf"
$ws.Range("O10").Value = "asdf"

$ws.Range("J11").Value = "asdf"
$ws.Range("O11").Value = "asdf"

$ws.Range("J12").Value = "asdfasdf"
$ws.Range("K12").Value = "asdf"
$ws.Range("L12").Value = "asdf"
$ws.Range("M12").Value = "asdfsa"
$ws.Range("N12").Value = "df"
$ws.Range("O12").Value = "asdf"

$ws.Range("K7").Select()
